$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.03133494016881249;  C = 0.618674247399176;   D = 0.7992451844515852;  E = 0.8940051367031318;  F = 0.9025266721007571;  G = 50 }
    3  = @{ B = 0.05905192170129663;  C = 0.7500407276805185;  D = 1.336363517888593;   E = 1.156011902139677;   F = 1.166466736141374;   G = 49 }
    4  = @{ B = 0.003554835396229139; C = 0.6262718728555511;  D = 0.7969898074296445;  E = 0.8927428562747756;  F = 0.9021829812466056;  G = 48 }
    5  = @{ B = -0.03511572363863584; C = 0.6635828519735995;  D = 0.8873739543468248;  E = 0.9420052836087623;  F = 0.9517525198146293;  G = 46 }
    6  = @{ B = -0.02357502177020739; C = 0.6102703326053934;  D = 0.7645395854549295;  E = 0.8743795431361198;  F = 0.8837201036899427;  G = 46 }
    7  = @{ B = 0.01163558867306234;  C = 0.5783626329962743;  D = 0.6903657101012227;  E = 0.8308824887438793;  F = 0.8432949623807621;  G = 34 }
    8  = @{ B = 0.009654851496814876; C = 0.5447572585674562;  D = 0.6494921422519612;  E = 0.8059107532797668;  F = 0.8183475071966614;  G = 33 }
    9  = @{ B = -0.04525758420911288; C = 0.5348296361048615;  D = 0.5246369517211453;  E = 0.7243182668697133;  F = 0.7466109761084062;  G = 16 }
    10 = @{ B = -0.120014643191906;   C = 0.3396452386262637;  D = 0.2564608965434484;  E = 0.5064196841982432;  F = 0.5186064467643261;  G = 10 }
    11 = @{ B = -0.248566666372188;   C = 0.301374533835348;   D = 0.1561541005302556;  E = 0.3951633846021865;  F = 0.3434543508584382 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
